# Auto-generated Excel COM-interop script that applies the
# "Updated cryptos list" scrape-refresh diff to the worksheet.
#
# All touched cells hold plain text (t="inlineStr" in the original
# file) even though several values look like numbers (e.g. "245.86").
# Assigning a bare string via .Value lets Excel's smart-typing turn
# those into real numbers, so for every cell we briefly force a Text
# number format, set the literal string, then ClearFormats() to drop
# the temporary format again (the cells carry no explicit style in
# the source workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '35.322.25'
Set-TextValue 'E2' '  +0.37%  '
Set-TextValue 'D3' '1.905.81'
Set-TextValue 'E3' '  +2.57%  '
Set-TextValue 'E4' '  -0.41%  '
Set-TextValue 'D5' '245.86'
Set-TextValue 'E5' '  +2.75%  '
Set-TextValue 'D6' '0.665'
Set-TextValue 'E6' '  +6.76%  '
Set-TextValue 'E7' '  -0.39%  '
Set-TextValue 'D8' '41.20'
Set-TextValue 'E8' '  -1.91%  '
Set-TextValue 'D9' '0.349'
Set-TextValue 'E9' '  +5.98%  '
Set-TextValue 'D10' '52.83'
Set-TextValue 'E10' '  +12.62%  '
Set-TextValue 'D11' '0.0718'
Set-TextValue 'E11' '  +3.55%  '
Set-TextValue 'D12' '0.0993'
Set-TextValue 'E12' '  +0.41%  '
Set-TextValue 'D13' '2.183.55'
Set-TextValue 'E13' '  +2.69%  '
Set-TextValue 'D14' '12.07'
Set-TextValue 'E14' '  +4.96%  '
Set-TextValue 'D15' '0.698'
Set-TextValue 'E15' '  +3.24%  '
Set-TextValue 'D16' '1.904.75'
Set-TextValue 'E16' '  +2.62%  '
Set-TextValue 'D17' '4.86'
Set-TextValue 'E17' '  +2.70%  '
Set-TextValue 'D18' '35.318.98'
Set-TextValue 'E18' '  +0.41%  '
Set-TextValue 'D19' '72.49'
Set-TextValue 'E19' '  +3.77%  '
Set-TextValue 'D20' '0.0₃0820'
Set-TextValue 'E20' '  +3.27%  '
Set-TextValue 'D21' '239.70'
Set-TextValue 'E21' '  -0.50%  '
Set-TextValue 'D22' '12.48'
Set-TextValue 'E22' '  +2.18%  '
Set-TextValue 'D23' '4.83'
Set-TextValue 'E23' '  +1.85%  '
Set-TextValue 'E24' '  -0.46%  '
Set-TextValue 'E25' '  +0.82%  '
Set-TextValue 'E26' '  +23.30%  '
Set-TextValue 'D27' '170.00'
Set-TextValue 'E27' '  +0.21%  '
Set-TextValue 'D28' '8.45'
Set-TextValue 'E28' '  +5.38%  '
Set-TextValue 'D29' '18.38'
Set-TextValue 'E29' '  +3.96%  '
Set-TextValue 'E30' '  +2.64%  '
Set-TextValue 'D31' '4.14'
Set-TextValue 'E31' '  +3.46%  '
Set-TextValue 'D32' '0.0565'
Set-TextValue 'E32' '  +0.61%  '
Set-TextValue 'B33' 'BinanceUSD'
Set-TextValue 'C33' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D33' '1.02'
Set-TextValue 'E33' '  +0.45%  '
Set-TextValue 'B34' 'ImmutableX'
Set-TextValue 'C34' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D34' '0.933'
Set-TextValue 'E34' '  +13.97%  '
Set-TextValue 'D35' '4.10'
Set-TextValue 'E35' '  +2.31%  '
Set-TextValue 'E36' '  -2.86%  '
Set-TextValue 'D37' '2.04'
Set-TextValue 'E37' '  -0.03%  '
Set-TextValue 'D39' '1.10'
Set-TextValue 'E39' '  +0.57%  '
Set-TextValue 'E40' '  +3.17%  '
Set-TextValue 'D41' '16.29'
Set-TextValue 'E41' '  +8.82%  '
Set-TextValue 'D42' '0.0637'
Set-TextValue 'E42' '  +6.08%  '
Set-TextValue 'D43' '89.85'
Set-TextValue 'E43' '  -0.13%  '
Set-TextValue 'D44' '1.339.86'
Set-TextValue 'E44' '  -0.61%  '
Set-TextValue 'D45' '2.38'
Set-TextValue 'E45' '  +2.60%  '
Set-TextValue 'D46' '47.51'
Set-TextValue 'E46' '  +37.01%  '
Set-TextValue 'B47' 'HuobiToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D47' '2.40'
Set-TextValue 'E47' '  -0.79%  '
Set-TextValue 'B48' 'MXToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D48' '2.78'
Set-TextValue 'E48' '  +1.52%  '
Set-TextValue 'E49' '  -0.31%  '
Set-TextValue 'D50' '2.091.76'
Set-TextValue 'E51' '  +3.44%  '
